$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NATALIA's balance (row 7, column C) from 7000 to 8000
$ws.Range("C7").Value = 8000

# Delete the MARCEL row (row 10) entirely, shifting rows below it up
$ws.Rows(10).Delete()
